$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.995.22'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.07%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.703.59'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.19'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3997'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4038'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.472'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '53.46'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.001'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08820'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '26.00'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.482'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.16%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001357'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.70%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.989'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.745.72'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +3.34%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.97'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07216'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.76'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.339'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.36'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.991.91'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.388'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.99%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.974'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.61'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +4.17%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +13.90%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '162.98'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '150.85'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +8.99%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.377'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.629'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +19.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.000.31'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +6.75%  '
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'VeChain'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.03175'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +6.33%  '
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08555'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.03%  '
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.258'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.45%  '
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.051'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.35%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2887'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.12'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.78%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09567'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8341'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.74%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.08'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.484'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.55%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.45%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.698'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.97%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7408'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.260'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.409'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08804'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +7.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.003'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '140.14'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.15%  '
